$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 1000
$ws.Range("F9").Value = 1472
$ws.Range("F11").Value = 1372
$ws.Range("F12").Value = 3026
$ws.Range("F13").Value = 482
$ws.Range("F14").Value = 1664
$ws.Range("F15").Value = 1367
$ws.Range("F17").Value = 249
$ws.Range("F18").Value = 1409
$ws.Range("F21").Value = 1142
$ws.Range("F22").Value = 11
$ws.Range("F23").Value = 408
$ws.Range("F25").Value = 3538
$ws.Range("F26").Value = 701
$ws.Range("F28").Value = 1568

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 48
$ws.Range("F4").Value = 173
$ws.Range("F9").Value = 30
$ws.Range("F10").Value = 15

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 798

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 798
$ws.Range("F7").Value = 48
$ws.Range("F8").Value = 173
$ws.Range("F14").Value = 30
$ws.Range("F15").Value = 15
$ws.Range("F17").Value = 1000
$ws.Range("F21").Value = 1472
$ws.Range("F23").Value = 1372
$ws.Range("F24").Value = 3026
$ws.Range("F25").Value = 482
$ws.Range("F26").Value = 1664
$ws.Range("F27").Value = 1367
$ws.Range("F29").Value = 249
$ws.Range("F30").Value = 1409
$ws.Range("F35").Value = 1142
$ws.Range("F36").Value = 11
$ws.Range("F37").Value = 408
$ws.Range("F39").Value = 3538
$ws.Range("F40").Value = 701
$ws.Range("F42").Value = 1568
